$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 2.88
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 3.2
$ws.Range("L3").Value = 3.5
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("Q3").Value = 2.08
$ws.Range("R3").Value = 1.73
$ws.Range("X3").Value = 12
$ws.Range("Y3").Value = 10
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 13
$ws.Range("AN3").Value = 4.5
$ws.Range("AO3").Value = 15
$ws.Range("AQ3").Value = 51
$ws.Range("AS3").Value = 201
$ws.Range("G5").Value = 2.55
$ws.Range("N5").Value = 10
$ws.Range("R5").Value = 1.75
$ws.Range("X5").Value = 12
$ws.Range("AD5").Value = 6
$ws.Range("AM5").Value = 34
$ws.Range("BC5").Value = 201
$ws.Range("G6").Value = 1.8
$ws.Range("H6").Value = 3.7
$ws.Range("J6").Value = 2.4
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 2.15
$ws.Range("R6").Value = 1.67
$ws.Range("U6").Value = 2.05
$ws.Range("V6").Value = 1.7
$ws.Range("W6").Value = 6
$ws.Range("Y6").Value = 9
$ws.Range("AA6").Value = 17
$ws.Range("AB6").Value = 34
$ws.Range("AC6").Value = 8.5
$ws.Range("AF6").Value = 67
$ws.Range("AG6").Value = 501
$ws.Range("AH6").Value = 10
$ws.Range("AI6").Value = 21
$ws.Range("AO6").Value = 9.5
$ws.Range("AP6").Value = 23
$ws.Range("AQ6").Value = 34
$ws.Range("AS6").Value = 201
$ws.Range("AU6").Value = 9
$ws.Range("AY6").Value = 29
$ws.Range("AZ6").Value = 41
